# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-holdings detail) right before the
#    "总计" (totals) summary sheet.
# 2) Insert a new row at the top of the "总计" sheet's data (just below the
#    header) summarizing the 2022-Q1 quarter, shifting the existing rows
#    down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# --- 1. Create the "2022-Q1" sheet -----------------------------------------
# Duplicate the "2021-Q4" sheet (same column layout/styling) and place the
# copy immediately before "总计", then rename it and replace its contents.
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$src.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Columns B-G hold text-formatted values (fund codes / decimal strings), so
# force Text number format before assigning values so the leading zeros in
# fund codes (e.g. "010695") are preserved instead of being parsed as numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "010695"
$newSheet.Range("C2").Value = "华夏磐益一年定期开放混合"
$newSheet.Range("D2").Value = "18.02"
$newSheet.Range("E2").Value = "82.41"
$newSheet.Range("F2").Value = "3.01"
$newSheet.Range("G2").Value = "0.5424"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006348"
$newSheet.Range("C3").Value = "银华盛利混合"
$newSheet.Range("D3").Value = "8.41"
$newSheet.Range("E3").Value = "88.43"
$newSheet.Range("F3").Value = "3.45"
$newSheet.Range("G3").Value = "0.2901"
$newSheet.Range("H3").Value = 4

# Drop the explicit Text format again now that the values are stored, so the
# cells end up with the same (default) style as the rest of the workbook.
$newSheet.Range("B2:G3").ClearFormats()

# --- 2. Update the "总计" sheet ---------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# The inserted row picks up the header's formatting; strip it back to the
# plain (unstyled) look used by the other data rows.
$ws.Range("B2:D2").ClearFormats()

# Column A keeps the bordered/centered "index" style (s=2) used throughout
# the sheet - copy it from the row below rather than re-deriving it.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.83

# Re-number the index column (A) for the rows that were shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# Restore the original active sheet/selection.
$wb.Worksheets.Item("2020-Q4").Activate()
